# Auto-generated from the xml diff: update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '26.245.55'
Set-TextCell $ws.Range("E2") '  +0.27%  '
Set-TextCell $ws.Range("D3") '1.593.82'
Set-TextCell $ws.Range("E4") '  -0.10%  '
Set-TextCell $ws.Range("D5") '211.80'
Set-TextCell $ws.Range("E5") '  +0.07%  '
Set-TextCell $ws.Range("D6") '0.504'
Set-TextCell $ws.Range("E6") '  -0.17%  '
Set-TextCell $ws.Range("E7") '  -0.06%  '
Set-TextCell $ws.Range("E8") '  -0.10%  '
Set-TextCell $ws.Range("D10") '18.91'
Set-TextCell $ws.Range("E10") '  -1.23%  '
Set-TextCell $ws.Range("E11") '  +0.67%  '
Set-TextCell $ws.Range("D12") '1.818.23'
Set-TextCell $ws.Range("E12") '  +0.54%  '
Set-TextCell $ws.Range("D13") '1.567.02'
Set-TextCell $ws.Range("E13") '  -1.15%  '
Set-TextCell $ws.Range("E14") '  -0.06%  '
Set-TextCell $ws.Range("D15") '0.503'
Set-TextCell $ws.Range("E15") '  -2.40%  '
Set-TextCell $ws.Range("D16") '63.56'
Set-TextCell $ws.Range("E16") '  -0.48%  '
Set-TextCell $ws.Range("D17") '26.213.91'
Set-TextCell $ws.Range("E17") '  +0.15%  '
Set-TextCell $ws.Range("D18") '229.22'
Set-TextCell $ws.Range("E18") '  +7.19%  '
Set-TextCell $ws.Range("D19") '0.0₃0720'
Set-TextCell $ws.Range("E19") '  -0.37%  '
Set-TextCell $ws.Range("D20") '7.59'
Set-TextCell $ws.Range("E20") '  +4.41%  '
Set-TextCell $ws.Range("D21") '0.999'
Set-TextCell $ws.Range("E21") '  -0.09%  '
Set-TextCell $ws.Range("D22") '4.23'
Set-TextCell $ws.Range("E22") '  -0.21%  '
Set-TextCell $ws.Range("B23") 'Avalanche'
Set-TextCell $ws.Range("C23") 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws.Range("D23") '8.89'
Set-TextCell $ws.Range("E23") '  -0.50%  '
Set-TextCell $ws.Range("B24") 'Toncoin'
Set-TextCell $ws.Range("C24") 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws.Range("D24") '2.13'
Set-TextCell $ws.Range("E24") '  +0.69%  '
Set-TextCell $ws.Range("D25") '145.66'
Set-TextCell $ws.Range("E25") '  +1.08%  '
Set-TextCell $ws.Range("E26") '  -0.13%  '
Set-TextCell $ws.Range("E27") '  +0.17%  '
Set-TextCell $ws.Range("D28") '0.113'
Set-TextCell $ws.Range("E28") '  +0.55%  '
Set-TextCell $ws.Range("E30") '  -0.49%  '
Set-TextCell $ws.Range("D31") '1.15'
Set-TextCell $ws.Range("E31") '  +0.30%  '
Set-TextCell $ws.Range("E32") '  +0.66%  '
Set-TextCell $ws.Range("D33") '1.458.97'
Set-TextCell $ws.Range("E33") '  +3.83%  '
Set-TextCell $ws.Range("E34") '  +0.36%  '
Set-TextCell $ws.Range("E36") '  +0.55%  '
Set-TextCell $ws.Range("E37") '  -3.99%  '
Set-TextCell $ws.Range("E38") '  -1.16%  '
Set-TextCell $ws.Range("D39") '0.819'
Set-TextCell $ws.Range("E39") '  -0.07%  '
Set-TextCell $ws.Range("D40") '5.77'
Set-TextCell $ws.Range("E40") '  -1.31%  '
Set-TextCell $ws.Range("E41") '  -0.10%  '
Set-TextCell $ws.Range("E42") '  +2.02%  '
Set-TextCell $ws.Range("D43") '0.929'
Set-TextCell $ws.Range("E43") '  -0.69%  '
Set-TextCell $ws.Range("D44") '1.730.74'
Set-TextCell $ws.Range("E44") '  +0.60%  '
Set-TextCell $ws.Range("D45") '0.756'
Set-TextCell $ws.Range("E45") '  -1.06%  '
Set-TextCell $ws.Range("D46") '60.40'
Set-TextCell $ws.Range("E46") '  -0.89%  '
Set-TextCell $ws.Range("D47") '87.56'
Set-TextCell $ws.Range("E47") '  +2.40%  '
Set-TextCell $ws.Range("E48") '  -0.18%  '
Set-TextCell $ws.Range("E49") '  +0.07%  '
Set-TextCell $ws.Range("E50") '  -0.16%  '
Set-TextCell $ws.Range("E51") '  -2.40%  '
